$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C on rows 2-21 changes from "SW" to "MS" (row 13's C cell was
# previously empty and now also gets "MS").
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 3).Value = "MS"
}

# Selection moved to C13 in the saved view state.
$ws.Range("C13").Select() | Out-Null
